# Regenerate save_data to use K (strikeouts) instead of Strike# (pitch count of strikes)
# in column G, for rows 2-12 of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 6
    4  = 7
    5  = 5
    6  = 8
    7  = 11
    8  = 6
    9  = 7
    10 = 2
    11 = 6
    12 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
